$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 66669290
$ws.Range("I64").Value = 200002020
$ws.Range("J64").Value = 2919.9
$ws.Range("K64").Value = 200002020
$ws.Range("L64").Value = 2919.9
$ws.Range("M64").Value = -200001772
$ws.Range("N64").Value = -3415.9

$ws.Range("H67").Value = 66669290
$ws.Range("I67").Value = 200002020
$ws.Range("J67").Value = 2919.9
$ws.Range("K67").Value = 200002020
$ws.Range("L67").Value = 2919.9
$ws.Range("M67").Value = -200001162
$ws.Range("N67").Value = -4635.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 35716810
$ws.Range("I61").Value = 50002332
$ws.Range("J61").Value = 3005.25
$ws.Range("K61").Value = 50002332
$ws.Range("L61").Value = 3005.25
$ws.Range("M61").Value = -50002120
$ws.Range("N61").Value = -3429.25

$ws.Range("H63").Value = 2661.6667
$ws.Range("I63").Value = 2403
$ws.Range("J63").Value = 2985
$ws.Range("K63").Value = 2403
$ws.Range("L63").Value = 2985
$ws.Range("M63").Value = -1717
$ws.Range("N63").Value = -4357

$ws.Range("H64").Value = 21400
$ws.Range("J64").Value = 21400
$ws.Range("L64").Value = 21400
$ws.Range("N64").Value = -21896

$ws.Range("H66").Value = 2661.6667
$ws.Range("I66").Value = 2403
$ws.Range("J66").Value = 2985
$ws.Range("K66").Value = 12015
$ws.Range("L66").Value = 14925
$ws.Range("M66").Value = -8583
$ws.Range("N66").Value = -21789

$ws.Range("H67").Value = 21400
$ws.Range("J67").Value = 21400
$ws.Range("L67").Value = 21400
$ws.Range("N67").Value = -23116

$ws.Range("H68").Value = 48000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 48000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H132").Value = 11186
$ws.Range("I132").Value = 11768
$ws.Range("J132").Value = 10652.5
$ws.Range("K132").Value = 35304
$ws.Range("L132").Value = 31957.5
$ws.Range("M132").Value = -32774
$ws.Range("N132").Value = -37017.5

$ws.Range("H136").Value = 35716810
$ws.Range("I136").Value = 50002332
$ws.Range("J136").Value = 3005.25
$ws.Range("K136").Value = 150006996
$ws.Range("L136").Value = 9015.75
$ws.Range("M136").Value = -150004446
$ws.Range("N136").Value = -14115.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3841.5737
$ws.Range("I31").Value = 1117.8948
$ws.Range("J31").Value = 8341.565000000001
$ws.Range("K31").Value = 1117.8948
$ws.Range("L31").Value = 8341.565000000001
$ws.Range("M31").Value = -822.8948
$ws.Range("N31").Value = -8931.565000000001

$ws.Range("H34").Value = 3841.5737
$ws.Range("I34").Value = 1117.8948
$ws.Range("J34").Value = 8341.565000000001
$ws.Range("K34").Value = 1117.8948
$ws.Range("L34").Value = 8341.565000000001
$ws.Range("M34").Value = -915.8948
$ws.Range("N34").Value = -8745.565000000001

$ws.Range("H99").Value = 1661.5714
$ws.Range("I99").Value = 1572
$ws.Range("J99").Value = 1822.8
$ws.Range("K99").Value = 1572
$ws.Range("L99").Value = 1822.8
$ws.Range("M99").Value = -74
$ws.Range("N99").Value = -4818.8

$ws.Range("H126").Value = 1661.5714
$ws.Range("I126").Value = 1572
$ws.Range("J126").Value = 1822.8
$ws.Range("K126").Value = 4716
$ws.Range("L126").Value = 5468.4
$ws.Range("M126").Value = -2246
$ws.Range("N126").Value = -10408.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 622.38336
$ws.Range("I5").Value = 439.7
$ws.Range("J5").Value = 987.75
$ws.Range("K5").Value = 1319.1
$ws.Range("L5").Value = 2963.25
$ws.Range("M5").Value = -1207.1
$ws.Range("N5").Value = -3187.25

$ws.Range("H115").Value = 988.5714
$ws.Range("I115").Value = 988.5714
$ws.Range("K115").Value = 2965.7142
$ws.Range("M115").Value = -1790.7142

$ws.Range("H122").Value = 797.7
$ws.Range("I122").Value = 442
$ws.Range("J122").Value = 1331.25
$ws.Range("K122").Value = 3978
$ws.Range("L122").Value = 11981.25
$ws.Range("M122").Value = -1528
$ws.Range("N122").Value = -16881.25

$ws.Range("H132").Value = 1054.1613
$ws.Range("I132").Value = 875.2222
$ws.Range("K132").Value = 7876.999800000001
$ws.Range("M132").Value = -5346.999800000001

$ws.Range("H135").Value = 622.38336
$ws.Range("I135").Value = 439.7
$ws.Range("J135").Value = 987.75
$ws.Range("K135").Value = 3957.3
$ws.Range("L135").Value = 8889.75
$ws.Range("M135").Value = -1422.3
$ws.Range("N135").Value = -13959.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3929.6428
$ws.Range("I80").Value = 2170
$ws.Range("J80").Value = 4409.5454
$ws.Range("K80").Value = 2170
$ws.Range("L80").Value = 4409.5454
$ws.Range("M80").Value = -1172
$ws.Range("N80").Value = -6405.5454

$ws.Range("H83").Value = 3929.6428
$ws.Range("I83").Value = 2170
$ws.Range("J83").Value = 4409.5454
$ws.Range("K83").Value = 10850
$ws.Range("L83").Value = 22047.727
$ws.Range("M83").Value = -5858
$ws.Range("N83").Value = -32031.727

$ws.Range("H122").Value = 1812
$ws.Range("I122").Value = 1815.4286
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5446.2858
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -2996.2858
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1527.1
$ws.Range("I68").Value = 1234.7273
$ws.Range("J68").Value = 1884.4445
$ws.Range("K68").Value = 1234.7273
$ws.Range("L68").Value = 1884.4445
$ws.Range("M68").Value = -485.7273
$ws.Range("N68").Value = -3382.4445

$ws.Range("H71").Value = 1527.1
$ws.Range("I71").Value = 1234.7273
$ws.Range("J71").Value = 1884.4445
$ws.Range("K71").Value = 6173.636500000001
$ws.Range("L71").Value = 9422.2225
$ws.Range("M71").Value = -2429.636500000001
$ws.Range("N71").Value = -16910.2225

$ws.Range("H93").Value = 1673.4
$ws.Range("I93").Value = 1387.7407
$ws.Range("J93").Value = 2637.5
$ws.Range("K93").Value = 1387.7407
$ws.Range("L93").Value = 2637.5
$ws.Range("M93").Value = -139.7407000000001
$ws.Range("N93").Value = -5133.5

$ws.Range("H140").Value = 35929
$ws.Range("J140").Value = 35929
$ws.Range("L140").Value = 35929
$ws.Range("N140").Value = -46289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3645.487
$ws.Range("I136").Value = 3939.3
$ws.Range("J136").Value = 2666.111
$ws.Range("K136").Value = 11817.9
$ws.Range("L136").Value = 7998.333
$ws.Range("M136").Value = -9267.900000000001
$ws.Range("N136").Value = -13098.333
